$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Several runs in this fixture contain a literal line-feed character (0x0A)
# embedded inside their <w:t> text (left over from how the fixture text was
# authored). This script strips/rewrites those embedded line feeds:
#   - at the end of a paragraph's only run, the line feed is simply dropped
#   - mid-paragraph (between two "plain" runs that would otherwise be
#     coalesced together by a same-formatting merge), the line feed is
#     replaced by a single space, and the run split is preserved by briefly
#     toggling formatting on the neighbouring run during the edit.
# ---------------------------------------------------------------------------

function Get-LineFeedRange([string]$needle) {
    # Locate the unique $needle (which itself contains a literal line feed)
    # and return a 1-character Range covering just the line feed.
    $found = $d.Content.Find.Execute($needle)
    if (-not $found) {
        return $null
    }
    $c = $d.Content
    $null = $c.Find.Execute($needle)
    $txt = $c.Text
    $idx = $txt.IndexOf([char]10)
    if ($idx -lt 0) {
        return $null
    }
    $lfStart = $c.Start + $idx
    return $d.Range($lfStart, $lfStart + 1)
}

function Remove-LineFeed([string]$needle) {
    $lfRange = Get-LineFeedRange $needle
    if ($null -eq $lfRange) {
        return $false
    }
    $lfRange.Text = ""
    return $true
}

function Replace-LineFeedWithSpace([string]$needle, [string]$guardText) {
    # $guardText is the text of the run immediately after the line feed.
    # Toggling its formatting around the edit stops the engine from
    # silently merging it into the edited (identically formatted) run.
    $guard1 = $d.Content
    $null = $guard1.Find.Execute($guardText)
    $guard1.Font.Bold = 1

    $lfRange = Get-LineFeedRange $needle
    if ($null -eq $lfRange) {
        $guard2 = $d.Content
        $null = $guard2.Find.Execute($guardText)
        $guard2.Font.Bold = 0
        return $false
    }
    $lfRange.Text = " "

    $guard2 = $d.Content
    $null = $guard2.Find.Execute($guardText)
    $guard2.Font.Bold = 0
    return $true
}

# "This is an annotatable resource in the casebook.\n" -> drop the trailing line feed
Remove-LineFeed "annotatable resource in the casebook.`n"

# ";\nreplaced: " -> "; replaced: " (line feed becomes a space; neighbouring
# runs already carry distinct character styles, so no guard is needed)
Replace-LineFeedWithSpace ";`nreplaced: " "foo bar baz"

# "; noted:\n" -> "; noted: " (line feed becomes a space; the following run
# "content to note;" shares identical formatting, so it needs the guard)
Replace-LineFeedWithSpace "; noted:`n" "content to note;"

# "This is the second chapter of the casebook.\n" -> drop the trailing line feed
Remove-LineFeed "chapter of the casebook.`n"
